# Fruta / hortaliza, semanal
#
# The daily-logic sheet is a rolling price log: each weekly refresh
# prepends one new record (row 145) and every previously logged record
# shifts down by one row, so the last existing record (old row 191)
# becomes row 192. The dimension grows from A1:T191 to A1:T192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 145..191 down by inserting a new blank row at 145.
$ws.Rows.Item(145).Insert()

# Populate the new row 145 with this week's record (same fixed
# attributes as the prior top record, new sample date and volume).
$ws.Cells.Item(145, 1).Value  = 5
$ws.Cells.Item(145, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(145, 3).Value  = "Maule"
$ws.Cells.Item(145, 4).Value  = 44559
$ws.Cells.Item(145, 5).Value  = 7
$ws.Cells.Item(145, 6).Value  = "Fruta"
$ws.Cells.Item(145, 7).Value  = 100108
$ws.Cells.Item(145, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(145, 9).Value  = 100108005
$ws.Cells.Item(145, 10).Value = "Piña"
$ws.Cells.Item(145, 11).Value = "Caramelo"
$ws.Cells.Item(145, 12).Value = "Segunda"
$ws.Cells.Item(145, 13).Value = 540
$ws.Cells.Item(145, 14).Value = 14000
$ws.Cells.Item(145, 15).Value = 14000
$ws.Cells.Item(145, 16).Value = 14000
$ws.Cells.Item(145, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(145, 18).Value = "Ecuador"
$ws.Cells.Item(145, 19).Value = 1000
$ws.Cells.Item(145, 20).Value = 14
